$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '42.447.31'
Set-TextValue 'E2' '  +2.66%  '
Set-TextValue 'D3' '2.235.17'
Set-TextValue 'E3' '  +2.16%  '
Set-TextValue 'E4' '  -0.15%  '
Set-TextValue 'D5' '252.06'
Set-TextValue 'E5' '  -1.06%  '
Set-TextValue 'D6' '0.629'
Set-TextValue 'E6' '  +0.47%  '
Set-TextValue 'D7' '69.36'
Set-TextValue 'E7' '  +1.48%  '
Set-TextValue 'E8' '  -0.09%  '
Set-TextValue 'D9' '0.630'
Set-TextValue 'E9' '  +7.69%  '
Set-TextValue 'D10' '40.03'
Set-TextValue 'E10' '  +5.72%  '
Set-TextValue 'D11' '59.79'
Set-TextValue 'E11' '  +1.66%  '
Set-TextValue 'D12' '0.0949'
Set-TextValue 'E12' '  +2.02%  '
Set-TextValue 'D13' '7.17'
Set-TextValue 'E13' '  -0.36%  '
Set-TextValue 'E14' '  -0.21%  '
Set-TextValue 'D15' '2.564.99'
Set-TextValue 'E15' '  +1.99%  '
Set-TextValue 'D16' '0.880'
Set-TextValue 'E16' '  +0.70%  '
Set-TextValue 'D17' '14.70'
Set-TextValue 'E17' '  +0.81%  '
Set-TextValue 'D18' '2.234.14'
Set-TextValue 'E18' '  +2.03%  '
Set-TextValue 'D19' '42.327.84'
Set-TextValue 'E19' '  +2.59%  '
Set-TextValue 'D20' '0.0₃0968'
Set-TextValue 'E20' '  +1.45%  '
Set-TextValue 'B21' 'Uniswap'
Set-TextValue 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D21' '6.21'
Set-TextValue 'E21' '  -0.05%  '
Set-TextValue 'B22' 'Litecoin'
Set-TextValue 'C22' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D22' '72.53'
Set-TextValue 'E22' '  +0.66%  '
Set-TextValue 'D23' '233.78'
Set-TextValue 'E23' '  +0.51%  '
Set-TextValue 'D24' '2.08'
Set-TextValue 'E24' '  +1.69%  '
Set-TextValue 'D25' '3.92'
Set-TextValue 'E25' '  +0.22%  '
Set-TextValue 'B26' 'Cosmos'
Set-TextValue 'C26' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D26' '11.46'
Set-TextValue 'E26' '  -3.34%  '
Set-TextValue 'B27' 'Dai'
Set-TextValue 'C27' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D27' '1.00'
Set-TextValue 'E27' '  +0.13%  '
Set-TextValue 'D28' '2.42'
Set-TextValue 'E28' '  -4.26%  '
Set-TextValue 'D29' '3.69'
Set-TextValue 'E29' '  -1.48%  '
Set-TextValue 'D30' '2.21'
Set-TextValue 'E30' '  +1.70%  '
Set-TextValue 'D31' '167.51'
Set-TextValue 'E31' '  -1.46%  '
Set-TextValue 'D32' '20.69'
Set-TextValue 'E32' '  +0.15%  '
Set-TextValue 'D33' '6.10'
Set-TextValue 'E33' '  +10.89%  '
Set-TextValue 'D34' '0.123'
Set-TextValue 'E34' '  +3.10%  '
Set-TextValue 'D35' '0.0785'
Set-TextValue 'E35' '  +5.36%  '
Set-TextValue 'D36' '0.124'
Set-TextValue 'E36' '  +0.17%  '
Set-TextValue 'D37' '28.17'
Set-TextValue 'E37' '  +3.78%  '
Set-TextValue 'D38' '4.68'
Set-TextValue 'E38' '  +1.27%  '
Set-TextValue 'D39' '4.13'
Set-TextValue 'E39' '  -1.92%  '
Set-TextValue 'D40' '0.0318'
Set-TextValue 'E40' '  +5.54%  '
Set-TextValue 'D41' '2.27'
Set-TextValue 'E41' '  +2.53%  '
Set-TextValue 'D42' '12.63'
Set-TextValue 'E42' '  -2.31%  '
Set-TextValue 'D43' '5.74'
Set-TextValue 'E43' '  +0.35%  '
Set-TextValue 'D44' '5.10'
Set-TextValue 'E44' '  +2.58%  '
Set-TextValue 'D45' '62.78'
Set-TextValue 'E45' '  -2.52%  '
Set-TextValue 'D46' '0.199'
Set-TextValue 'E46' '  -1.59%  '
Set-TextValue 'D47' '8.68'
Set-TextValue 'E47' '  +0.55%  '
Set-TextValue 'D48' '0.101'
Set-TextValue 'E48' '  -0.13%  '
Set-TextValue 'E49' '  -0.27%  '
Set-TextValue 'D50' '1.17'
Set-TextValue 'E50' '  +2.78%  '
Set-TextValue 'B51' 'TrustWalletToken'
Set-TextValue 'C51' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D51' '1.18'
Set-TextValue 'E51' '  +0.08%  '
